# Update the dSF column (F) values for the affected rows, per the repull
# of data / mean calculation described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = 4
$ws.Range("F10").Value = -6
